$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 6, shifting rows 6:45 down to 7:46
$ws.Rows.Item(6).Insert()

# Copy the row's fixed/common values and formatting from the row above (row 5)
$ws.Rows.Item(5).Copy()
$ws.Rows.Item(6).PasteSpecial()

# Fill in the values for the new row 6
$ws.Cells.Item(6, 1).Value = 10
$ws.Cells.Item(6, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(6, 3).Value = 'La Araucanía'
$ws.Cells.Item(6, 4).Value = 45149
$ws.Cells.Item(6, 5).Value = 9
$ws.Cells.Item(6, 6).Value = 'Fruta'
$ws.Cells.Item(6, 7).Value = 100108
$ws.Cells.Item(6, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(6, 9).Value = 100108001
$ws.Cells.Item(6, 10).Value = 'Guayaba'
$ws.Cells.Item(6, 11).Value = 'Sin especificar'
$ws.Cells.Item(6, 12).Value = 'Primera'
$ws.Cells.Item(6, 13).Value = 100
$ws.Cells.Item(6, 14).Value = 2700
$ws.Cells.Item(6, 15).Value = 2700
$ws.Cells.Item(6, 16).Value = 2700
$ws.Cells.Item(6, 17).Value = '$/kilo'
$ws.Cells.Item(6, 18).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(6, 19).Value = 2700
$ws.Cells.Item(6, 20).Value = 1
